# Updates betting-odds base rows for "Mexico Liga MX" sheet.
# The underlying source rows got re-sorted/re-matched against their correct
# match ids, which results in whole data-rows (every column except the
# positional index in column A) being swapped between specific row pairs
# (and one 3-way rotation), plus a direct odds refresh for the still-open
# fixture in the final row.

function Set-RowData {
    param($Row, $B, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U, $V, $W, $X, $Y, $Z, $AA, $AB)
    $ws = $wb.ActiveSheet
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
    $ws.Cells.Item($Row, 21).Value = $U
    $ws.Cells.Item($Row, 22).Value = $V
    $ws.Cells.Item($Row, 23).Value = $W
    $ws.Cells.Item($Row, 24).Value = $X
    $ws.Cells.Item($Row, 25).Value = $Y
    $ws.Cells.Item($Row, 26).Value = $Z
    $ws.Cells.Item($Row, 27).Value = $AA
    $ws.Cells.Item($Row, 28).Value = $AB
}

$wb = $excel.ActiveWorkbook

Set-RowData 36 6754017 "Chivas Guadalajara" "Tijuana" 1 0 "H" 1.5 3 1.5 1.615 4.2 5 -0.75 1.8 2.05 2.75 2 1.85 0.615 -1 -1 0.4 -0.5 -1 0.8500000000000001
Set-RowData 37 6754019 "Mazatlan FC" "Puebla" 1 0 "H" 2.4 3.3 2.9 1.85 3.75 4 -0.5 1.875 1.975 3 1.825 2.025 0.8500000000000001 -1 -1 0.875 -1 -1 1.025

Set-RowData 72 6754049 "Juarez FC" "Atlas" 1 2 "A" 2.75 3.25 2.375 2.6 3.2 2.8 0 1.85 2 2.25 2.1 1.775 -1 -1 1.8 -1 1 1.1 -1
Set-RowData 73 6754048 "Atletico San Luis" "Mazatlan FC" 3 2 "H" 1.615 4 4.5 1.6 4.5 5 -1 1.95 1.9 3 1.925 1.925 0.6000000000000001 -1 -1 0 0 0.925 -1

Set-RowData 94 6754067 "Tigres UANL" "Toluca" 2 2 "D" 1.8 3.3 4.333 1.533 4.2 6 -1 1.925 1.925 3 2 1.85 -1 3.2 -1 -1 0.925 1 -1
Set-RowData 95 7260442 "Santos Laguna" "Tijuana" 2 1 "H" 1.75 3.6 4.2 1.65 4 4.75 -0.75 1.8 2.05 3 1.85 2 0.6499999999999999 -1 -1 0.4 -0.5 0 0

Set-RowData 98 6754074 "Chivas Guadalajara" "Atlas" 4 1 "H" 2.3 3.3 2.8 2.4 3.2 3.1 -0.25 2.15 1.725 2 1.925 1.925 1.4 -1 -1 1.15 -1 0.925 -1
Set-RowData 99 6754641 "Pachuca" "Tigres UANL" 1 1 "D" 2.875 3.5 2.2 2.9 3.5 2.4 0.25 1.775 2.1 2.5 1.825 2.025 -1 2.5 -1 0.3875 -0.5 -1 1.025

Set-RowData 126 6754096 "Queretaro" "Chivas Guadalajara" 1 2 "A" 3.25 3.3 2.2 2.8 3.1 2.7 0 2 1.85 2.25 1.925 1.925 -1 -1 1.7 -1 0.8500000000000001 0.925 -1
Set-RowData 127 6754097 "Toluca" "Puebla" 0 1 "A" 1.5 4 7 1.45 4.5 7 -1.25 2 1.85 3.25 1.975 1.875 -1 -1 6 -1 0.8500000000000001 -1 0.875

Set-RowData 130 6754100 "Atlas" "Pachuca" 0 2 "A" 2.1 3.3 3.5 2.45 3.2 3 -0.25 2.05 1.8 2.25 1.875 1.975 -1 -1 2 -1 0.8 -0.5 0.4875
Set-RowData 131 6754101 "Cruz Azul" "Juarez FC" 2 0 "H" 2.25 3.3 3.25 1.65 3.8 5.25 -0.75 1.8 2.05 2.75 1.925 1.925 0.6499999999999999 -1 -1 0.8 -1 -1 0.925

Set-RowData 132 6754103 "Mazatlan FC" "Santos Laguna" 3 1 "H" 2.5 3.3 2.75 2.8 3.4 2.5 0 2.05 1.8 3 1.975 1.875 1.8 -1 -1 1.05 -1 0.9750000000000001 -1
Set-RowData 133 6754129 "Atletico San Luis" "Club America" 0 1 "A" 3.75 3.3 2 5 3.8 1.7 0.75 1.95 1.9 3 1.925 1.925 -1 -1 0.7 -0.5 0.45 -1 0.925

Set-RowData 138 6754643 "Pachuca" "Monterrey" 0 2 "A" 2.7 3.4 2.5 2.7 3.5 2.55 0 1.975 1.875 2.75 2.025 1.825 -1 -1 1.55 -1 0.875 -1 0.825
Set-RowData 139 6754105 "Chivas Guadalajara" "Cruz Azul" 1 0 "H" 1.8 3.75 4.2 2.45 3.25 3 -0.25 2.05 1.75 2.5 1.975 1.875 1.45 -1 -1 1.05 -1 -1 0.875

Set-RowData 175 7612675 "Monterrey" "Puebla" 2 0 "H" 1.571 3.75 6 1.363 5 7.5 -1.5 2.025 1.825 3 2 1.85 0.363 -1 -1 1.025 -1 -1 0.8500000000000001
Set-RowData 176 7612685 "Tijuana" "Club America" 0 2 "A" 4.5 3.5 1.8 2.45 3.3 2.875 -0.25 2.1 1.775 2.25 1.8 2.05 -1 -1 1.875 -1 0.7749999999999999 -0.5 0.5249999999999999

Set-RowData 200 7713694 "Cruz Azul" "Tijuana" 1 0 "H" 1.727 3.75 4.5 1.533 4.2 6 -1 1.95 1.9 2.75 1.95 1.9 0.5329999999999999 -1 -1 0 0 -1 0.8999999999999999
Set-RowData 201 7612821 "Mazatlan FC" "Leon" 2 2 "D" 2.5 3.3 2.75 3.3 3.6 2.15 0.25 1.925 1.925 2.5 1.875 1.975 -1 2.6 -1 0.4625 -0.5 0.875 -1

Set-RowData 264 7612885 "Santos Laguna" "Cruz Azul" 3 0 "H" 3.3 3.6 2.05 4.5 4 1.727 0.75 1.875 1.975 2.75 1.975 1.875 3.5 -1 -1 0.875 -1 0.4875 -0.5
Set-RowData 265 7612884 "Pachuca" "Queretaro" 1 2 "A" 1.6 4.2 5 2.05 3.8 3.3 -0.5 2.05 1.8 2.75 1.925 1.925 -1 -1 2.3 -1 0.8 0.4625 -0.5

Set-RowData 303 7745553 "Unam Pumas" "Leon" 1 0 "H" 2.2 3.5 3 1.909 3.8 3.8 -0.5 1.975 1.875 3 1.8 2.05 0.909 -1 -1 0.9750000000000001 -1 -1 1.05
Set-RowData 304 7745552 "Atlas" "Atletico San Luis" 2 1 "H" 1.833 3.6 4.2 2.375 3.4 3 -0.25 2.05 1.8 2.75 1.85 2 1.375 -1 -1 1.05 -1 0.425 -0.5

Set-RowData 310 7612930 "Chivas Guadalajara" "Queretaro" 2 0 "H" 2.4 3.2 2.75 1.533 4.2 6 -1 1.9 1.95 2.75 1.975 1.875 0.5329999999999999 -1 -1 0.8999999999999999 -1 -1 0.875
Set-RowData 311 7612929 "Tigres UANL" "Necaxa" 5 2 "H" 2.25 3.25 2.9 1.533 4.333 5.5 -1 1.9 1.95 2.75 1.85 2 0.5329999999999999 -1 -1 0.8999999999999999 -1 0.8500000000000001 -1
Set-RowData 312 7612933 "Santos Laguna" "Pachuca" 0 2 "A" 2.5 3.2 2.625 2.625 3.4 2.6 0 1.975 1.875 2.75 1.9 1.95 -1 -1 1.6 -1 0.875 -1 0.95

# Row 339 is the still-unplayed fixture (no FTHG/FTAG/FTR yet); only its
# odds columns M..U were refreshed.
$ws339 = $wb.ActiveSheet
$ws339.Cells.Item(339, 13).Value = 3.1
$ws339.Cells.Item(339, 14).Value = 3.25
$ws339.Cells.Item(339, 15).Value = 2.4
$ws339.Cells.Item(339, 16).Value = 0.25
$ws339.Cells.Item(339, 17).Value = 1.8
$ws339.Cells.Item(339, 18).Value = 2.05
$ws339.Cells.Item(339, 19).Value = 2.25
$ws339.Cells.Item(339, 20).Value = 1.975
$ws339.Cells.Item(339, 21).Value = 1.875
